$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "487×2=974" "390×3=1170"
Replace-Text "493×9=4437" "780×9=7020"
Replace-Text "918×4=3672" "789×6=4734"
Replace-Text "496×5=2480" "844×3=2532"
Replace-Text "368×3=1104" "398×7=2786"
Replace-Text "360×8=2880" "927×4=3708"
Replace-Text "903×8=7224" "947×3=2841"
Replace-Text "744×9=6696" "821×4=3284"
Replace-Text "972×6=5832" "302×4=1208"
Replace-Text "886×8=7088" "128×4=512"
Replace-Text "793×8=6344" "775×8=6200"
Replace-Text "623×2=1246" "381×2=762"
Replace-Text "633×6=3798" "896×7=6272"
Replace-Text "969×5=4845" "268×8=2144"
Replace-Text "380×2=760" "116×8=928"
Replace-Text "439×2=878" "499×9=4491"
Replace-Text "787×2=1574" "132×2=264"
Replace-Text "900×6=5400" "762×6=4572"
Replace-Text "751×5=3755" "764×3=2292"
Replace-Text "432×4=1728" "225×5=1125"
Replace-Text "992×9=8928" "221×6=1326"
Replace-Text "361×4=1444" "574×9=5166"
Replace-Text "807×3=2421" "697×2=1394"
Replace-Text "484×8=3872" "930×9=8370"
Replace-Text "283×8=2264" "760×5=3800"
